# Updated cryptos list on Wed Oct 23 19:18:14 UTC 2024 with GitHub Actions
# Refresh the per-coin Price (column D) and Volume(1h) (column E) figures that
# the scraper pulled this run; rows 49/50 also swap rank (Optimism now ranks
# above BabyDogeCoin) so their Coin/Link/Price/Volume cells trade places.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "65.891.89"
$ws.Range("E2").Value = "  -2.36%  "

# Row 3
$ws.Range("D3").Value = "2.490.55"
$ws.Range("E3").Value = "  -5.26%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.06%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.31"
$ws.Range("E5").Value = "  -2.67%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.14"
$ws.Range("E6").Value = "  +0.30%  "

# Row 7
$ws.Range("E7").Value = "  +0.13%  "

# Row 8
$ws.Range("E8").Value = "  -3.03%  "

# Row 9
$ws.Range("D9").Value = "2.491.14"
$ws.Range("E9").Value = "  -5.21%  "

# Row 10
$ws.Range("E10").Value = "  -3.07%  "

# Row 11
$ws.Range("E11").Value = "  -0.42%  "

# Row 12
$ws.Range("E12").Value = "  -4.93%  "

# Row 13
$ws.Range("E13").Value = "  -2.79%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.26"
$ws.Range("E14").Value = "  -5.11%  "

# Row 15
$ws.Range("D15").Value = "2.939.51"
$ws.Range("E15").Value = "  -5.34%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000173"
$ws.Range("E16").Value = "  -5.03%  "

# Row 17
$ws.Range("D17").Value = "65.724.76"
$ws.Range("E17").Value = "  -2.12%  "

# Row 18
$ws.Range("D18").Value = "2.499.32"
$ws.Range("E18").Value = "  -4.91%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.12"
$ws.Range("E19").Value = "  -7.68%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.59"
$ws.Range("E20").Value = "  -5.39%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "342.40"
$ws.Range("E21").Value = "  -4.04%  "

# Row 22
$ws.Range("E22").Value = "  -3.69%  "

# Row 23
$ws.Range("E23").Value = "  -3.15%  "

# Row 24
$ws.Range("E24").Value = "  -0.08%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.93"
$ws.Range("E25").Value = "  -0.06%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "68.77"
$ws.Range("E26").Value = "  -1.19%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.80"
$ws.Range("E27").Value = "  -5.10%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  +0.00%  "

# Row 29
$ws.Range("D29").Value = "2.618.94"
$ws.Range("E29").Value = "  -5.16%  "

# Row 30
$ws.Range("D30").Value = "0.0₃0956"
$ws.Range("E30").Value = "  -5.12%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "521.06"
$ws.Range("E31").Value = "  -4.42%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.00"
$ws.Range("E32").Value = "  +1.02%  "

# Row 33
$ws.Range("E33").Value = "  -3.04%  "

# Row 34
$ws.Range("E34").Value = "  -4.48%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.130"
$ws.Range("E35").Value = "  -3.83%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.01%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "156.28"
$ws.Range("E37").Value = "  -0.06%  "

# Row 38
$ws.Range("E38").Value = "  -4.26%  "

# Row 39
$ws.Range("E39").Value = "  -3.16%  "

# Row 40
$ws.Range("E40").Value = "  +0.62%  "

# Row 41
$ws.Range("E41").Value = "  -4.24%  "

# Row 42
$ws.Range("E42").Value = "  -3.26%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.03"
$ws.Range("E43").Value = "  -3.63%  "

# Row 44
$ws.Range("E44").Value = "  -0.04%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.41"
$ws.Range("E45").Value = "  -0.55%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "145.85"
$ws.Range("E46").Value = "  -4.57%  "

# Row 47
$ws.Range("E47").Value = "  -4.92%  "

# Row 48
$ws.Range("E48").Value = "  -3.99%  "

# Row 49
$ws.Range("B49").Value = "Optimism"
$ws.Range("C49").Value = "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.70"
$ws.Range("E49").Value = "  +0.49%  "

# Row 50
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₆0268"
$ws.Range("E50").Value = "  -9.54%  "

# Row 51
$ws.Range("E51").Value = "  -2.78%  "
